# Applies the error_report.xlsx edit: rows 2-13 of Sheet1 updated/replaced
# per the new validation run (75k-row edition). Row 13 is newly added.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 93
$ws.Cells.Item(2, 2).Value = "invoice_number"
$ws.Cells.Item(2, 3).Value = "inv-91"
$ws.Cells.Item(2, 4).Value = "COMPOSITE_DUPLICATE"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = "Original row flagged: Conflict with row 447"

# Row 3
$ws.Cells.Item(3, 1).Value = 102
$ws.Cells.Item(3, 2).Value = "customer_code"
$ws.Cells.Item(3, 3).Value = "'"
$ws.Cells.Item(3, 4).Value = "MISSING_REQUIRED"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = "Mandatory field is empty"

# Row 4
$ws.Cells.Item(4, 1).Value = 402
$ws.Cells.Item(4, 2).Value = "customer_code"
$ws.Cells.Item(4, 3).Value = "CUST-X"
$ws.Cells.Item(4, 4).Value = "PATTERN_MISMATCH"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = "Invalid pattern: CUST-X"

# Row 5
$ws.Cells.Item(5, 1).Value = 403
$ws.Cells.Item(5, 2).Value = "customer_code"
$ws.Cells.Item(5, 3).Value = "CUST-X"
$ws.Cells.Item(5, 4).Value = "PATTERN_MISMATCH"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = "Invalid pattern: CUST-X"

# Row 6
$ws.Cells.Item(6, 1).Value = 447
$ws.Cells.Item(6, 2).Value = "invoice_number"
$ws.Cells.Item(6, 3).Value = "INV-91"
$ws.Cells.Item(6, 4).Value = "COMPOSITE_DUPLICATE"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = "Duplicate of row 93"

# Row 7
$ws.Cells.Item(7, 1).Value = 503
$ws.Cells.Item(7, 2).Value = "customer_code"
$ws.Cells.Item(7, 3).Value = "CUST-Y"
$ws.Cells.Item(7, 4).Value = "PATTERN_MISMATCH"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = "Invalid pattern: CUST-Y"

# Row 8
$ws.Cells.Item(8, 1).Value = 504
$ws.Cells.Item(8, 2).Value = "customer_code"
$ws.Cells.Item(8, 3).Value = "CUST-Y"
$ws.Cells.Item(8, 4).Value = "PATTERN_MISMATCH"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = "Invalid pattern: CUST-Y"

# Row 9
$ws.Cells.Item(9, 1).Value = 505
$ws.Cells.Item(9, 2).Value = "customer_code"
$ws.Cells.Item(9, 3).Value = "CUST-Y"
$ws.Cells.Item(9, 4).Value = "PATTERN_MISMATCH"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = "Invalid pattern: CUST-Y"

# Row 10
$ws.Cells.Item(10, 1).Value = 30
$ws.Cells.Item(10, 2).Value = "name"
$ws.Cells.Item(10, 3).Value = "CompanyCompanyCompanyCompanyCompanyCompanyCompanyCompanyCompanyCompanyCompanyCompanyCompanyCompanyCompanyCompanyCompanyCompany 28"
$ws.Cells.Item(10, 4).Value = "BUSINESS_RULE_VIOLATION"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = "Too long (Max: 50)"

# Row 11
$ws.Cells.Item(11, 1).Value = 302
$ws.Cells.Item(11, 2).Value = "amount"
$ws.Cells.Item(11, 3).Value = "'5000"
$ws.Cells.Item(11, 4).Value = "BUSINESS_RULE_VIOLATION"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = "Constraint: Since status is suspend, this must be None"

# Row 12
$ws.Cells.Item(12, 1).Value = 703
$ws.Cells.Item(12, 2).Value = "email"
$ws.Cells.Item(12, 3).Value = "invalid_email.com"
$ws.Cells.Item(12, 4).Value = "PATTERN_MISMATCH"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = "Invalid email format"

# Row 13
$ws.Cells.Item(13, 1).Value = 202
$ws.Cells.Item(13, 2).Value = "project_code"
$ws.Cells.Item(13, 3).Value = "WRONG-CODE"
$ws.Cells.Item(13, 4).Value = "PATTERN_MISMATCH"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = "Invalid pattern: WRONG-CODE"
